$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Overall Status cell (C22) changes from "Pass" to "Fail"
$ws.Range("C22").Value = "Fail"

# Update the active selection to C23 (as reflected in the saved view state)
$ws.Activate()
$ws.Range("C23").Select()
